# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# The "Estado de Cuenta" detail table (rows 16-80, columns C:F) was
# re-sorted: previously grouped by worker (5 workers x 13 periods each,
# periods descending 2102..2002 per worker), now grouped by period
# ascending (2002..2102), with the 5 workers repeated for each period.
# Column F (Valor Mora) stays 35112 for every period except the most
# recent one (2102), which keeps 25749.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$workers = @(
    @('9153983',    'JUAN MANUEL LOMBANA DIAZ'),
    @('1047460015', 'LAURA PATRICIA MUÑOZ CARABALLO'),
    @('45555972',   'OLGA PATRICIA POMBO SOTO'),
    @('80874921',   'SAMUEL JOSE POMBO COGOLLO'),
    @('1007130691', 'JESUS DAVID PEREZ MARTINEZ')
)

$periods = @('2002','2003','2004','2005','2006','2007','2008','2009','2010','2011','2012','2101','2102')

$row = 16
foreach ($period in $periods) {
    foreach ($worker in $workers) {
        if ($period -eq '2102') {
            $valorMora = 25749
        } else {
            $valorMora = 35112
        }

        $ws.Range("C$row").Value = $worker[0]
        $ws.Range("D$row").Value = $worker[1]
        $ws.Range("E$row").Value = $period
        $ws.Range("F$row").Value = $valorMora

        $row = $row + 1
    }
}
